# Generate Report for Handoff
#
# Re-sorts the per-file rows (alphabetically by source file name) on the
# Overview / zh-cn / de-de sheets and marks the 604ee6db-... file as
# "Ready for handoff" (with a refreshed Latest Handoff Datetime), since it
# is now queued for a new handoff instead of being in sync.

$wb = $excel.ActiveWorkbook

function Set-CellAndHyperlink {
    param(
        $ws,
        [string]$cellRef,
        [string]$newValue
    )

    $ws.Range($cellRef).Value = $newValue

    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq ("`$" + $cellRef.Substring(0,1) + "`$" + $cellRef.Substring(1))) {
            $hl.TextToDisplay = $newValue
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

Set-CellAndHyperlink $ov "A2" "fffffa81bd4d-fc1e-4e68-8dfc-62ab0f1c43fa.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"

Set-CellAndHyperlink $ov "A3" "ffffff2ff9b7f8-3bbc-4723-a2df-446ecb18758e.md"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"

Set-CellAndHyperlink $ov "A4" "604ee6db-e269-4e9f-9b7d-bb6c8fc97582.md"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"

# Row 5 (.localization-config) is unchanged.

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

Set-CellAndHyperlink $zh "A2" "fffffa81bd4d-fc1e-4e68-8dfc-62ab0f1c43fa.md"
$zh.Range("B2").Value = "Handed back: in sync with en-US"
Set-CellAndHyperlink $zh "C2" "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.zh-cn.xlf"
$zh.Range("D2").Value = "2016-03-09 12:58:15"
Set-CellAndHyperlink $zh "E2" "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.md"
Set-CellAndHyperlink $zh "F2" "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.zh-cn.xlf"
$zh.Range("G2").Value = "2016-03-09 12:58:49"
$zh.Range("H2").Value = "Include"

Set-CellAndHyperlink $zh "A3" "ffffff2ff9b7f8-3bbc-4723-a2df-446ecb18758e.md"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
Set-CellAndHyperlink $zh "C3" "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.zh-cn.xlf"
$zh.Range("D3").Value = "2016-03-09 12:58:15"
Set-CellAndHyperlink $zh "E3" "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.md"
Set-CellAndHyperlink $zh "F3" "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.zh-cn.xlf"
$zh.Range("G3").Value = "2016-03-09 12:58:49"
$zh.Range("H3").Value = "Include"

Set-CellAndHyperlink $zh "A4" "604ee6db-e269-4e9f-9b7d-bb6c8fc97582.md"
$zh.Range("B4").Value = "Ready for handoff"
Set-CellAndHyperlink $zh "C4" "604ee6db-e269-4e9f-9b7d-bb6c8fc97582.d3c932d3f41c953547e0b93e6c368f73f27891d4.zh-cn.xlf"
$zh.Range("D4").Value = "2016-03-09 13:01:29"
Set-CellAndHyperlink $zh "E4" "604ee6db-e269-4e9f-9b7d-bb6c8fc97582.md"
Set-CellAndHyperlink $zh "F4" "604ee6db-e269-4e9f-9b7d-bb6c8fc97582.d3c932d3f41c953547e0b93e6c368f73f27891d4.zh-cn.xlf"
$zh.Range("G4").Value = "2016-03-09 13:00:38"
$zh.Range("H4").Value = "Include"

# Row 5 (.localization-config) is unchanged.

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

Set-CellAndHyperlink $de "A2" "fffffa81bd4d-fc1e-4e68-8dfc-62ab0f1c43fa.md"
$de.Range("B2").Value = "Handed back: in sync with en-US"
Set-CellAndHyperlink $de "C2" "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.de-de.xlf"
$de.Range("D2").Value = "2016-03-09 12:58:22"
Set-CellAndHyperlink $de "E2" "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.md"
Set-CellAndHyperlink $de "F2" "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.de-de.xlf"
$de.Range("G2").Value = "2016-03-09 12:59:05"
$de.Range("H2").Value = "Include"

Set-CellAndHyperlink $de "A3" "ffffff2ff9b7f8-3bbc-4723-a2df-446ecb18758e.md"
$de.Range("B3").Value = "Handed back: in sync with en-US"
Set-CellAndHyperlink $de "C3" "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.de-de.xlf"
$de.Range("D3").Value = "2016-03-09 12:58:22"
Set-CellAndHyperlink $de "E3" "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.md"
Set-CellAndHyperlink $de "F3" "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.de-de.xlf"
$de.Range("G3").Value = "2016-03-09 12:59:05"
$de.Range("H3").Value = "Include"

Set-CellAndHyperlink $de "A4" "604ee6db-e269-4e9f-9b7d-bb6c8fc97582.md"
$de.Range("B4").Value = "Ready for handoff"
Set-CellAndHyperlink $de "C4" "604ee6db-e269-4e9f-9b7d-bb6c8fc97582.d3c932d3f41c953547e0b93e6c368f73f27891d4.de-de.xlf"
$de.Range("D4").Value = "2016-03-09 13:01:41"
Set-CellAndHyperlink $de "E4" "604ee6db-e269-4e9f-9b7d-bb6c8fc97582.md"
Set-CellAndHyperlink $de "F4" "604ee6db-e269-4e9f-9b7d-bb6c8fc97582.d3c932d3f41c953547e0b93e6c368f73f27891d4.de-de.xlf"
$de.Range("G4").Value = "2016-03-09 13:00:54"
$de.Range("H4").Value = "Include"

# Row 5 (.localization-config) is unchanged.
